# Applies the "Updated cryptos list" data refresh to Sheet1 (columns D = Price, E = Volume(1h)).
# Rows 14/15 additionally swap coin identity (Toncoin <-> WrappedliquidstakedEther2.0).
# Numeric-looking Price values are prefixed with a literal leading apostrophe so Excel
# keeps them as text (matching the original inline-string cell type) instead of parsing
# them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '86.593.60'
$ws.Range("E2").Value = '  -3.49%  '

$ws.Range("D3").Value = '3.133.68'
$ws.Range("E3").Value = '  -7.51%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''203.57'
$ws.Range("E5").Value = '  -7.78%  '

$ws.Range("D6").Value = '''605.20'
$ws.Range("E6").Value = '  -7.10%  '

$ws.Range("E7").Value = '  -9.21%  '

$ws.Range("E8").Value = '  +7.02%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '3.133.81'
$ws.Range("E10").Value = '  -7.16%  '

$ws.Range("D11").Value = '''0.524'
$ws.Range("E11").Value = '  -12.15%  '

$ws.Range("E12").Value = '  +4.36%  '

$ws.Range("E13").Value = '  -17.63%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.712.64'
$ws.Range("E14").Value = '  -7.07%  '

$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '''5.19'
$ws.Range("E15").Value = '  -6.56%  '

$ws.Range("D16").Value = '86.320.02'
$ws.Range("E16").Value = '  -3.67%  '

$ws.Range("D17").Value = '''31.68'
$ws.Range("E17").Value = '  -14.60%  '

$ws.Range("D18").Value = '3.124.87'
$ws.Range("E18").Value = '  -7.10%  '

$ws.Range("D19").Value = '''2.93'
$ws.Range("E19").Value = '  -7.95%  '

$ws.Range("D20").Value = '''13.20'
$ws.Range("E20").Value = '  -11.27%  '

$ws.Range("D21").Value = '''407.99'
$ws.Range("E21").Value = '  -11.46%  '

$ws.Range("E22").Value = '  -13.13%  '

$ws.Range("D23").Value = '''5.02'
$ws.Range("E23").Value = '  -9.27%  '

$ws.Range("D24").Value = '''5.06'
$ws.Range("E24").Value = '  -9.90%  '

$ws.Range("D25").Value = '''11.58'
$ws.Range("E25").Value = '  -10.77%  '

$ws.Range("D26").Value = '3.299.12'

$ws.Range("D27").Value = '''72.43'
$ws.Range("E27").Value = '  -8.21%  '

$ws.Range("E28").Value = '  -10.74%  '

$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("D30").Value = '''0.160'
$ws.Range("E30").Value = '  -22.16%  '

$ws.Range("E31").Value = '  -0.53%  '

$ws.Range("D32").Value = '''528.89'
$ws.Range("E32").Value = '  -11.70%  '

$ws.Range("D33").Value = '''8.17'
$ws.Range("E33").Value = '  -12.84%  '

$ws.Range("D34").Value = '''1.83'
$ws.Range("E34").Value = '  -13.85%  '

$ws.Range("D35").Value = '''1.27'
$ws.Range("E35").Value = '  -20.66%  '

$ws.Range("D36").Value = '''6.47'
$ws.Range("E36").Value = '  -12.85%  '

$ws.Range("D37").Value = '''0.131'
$ws.Range("E37").Value = '  -9.12%  '

$ws.Range("D38").Value = '''21.40'
$ws.Range("E38").Value = '  -8.81%  '

$ws.Range("D39").Value = '''0.999'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("E41").Value = '  -7.82%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  -14.65%  '

$ws.Range("D44").Value = '''0.364'
$ws.Range("E44").Value = '  -15.16%  '

$ws.Range("D45").Value = '''148.83'
$ws.Range("E45").Value = '  -5.57%  '

$ws.Range("D46").Value = '''169.90'
$ws.Range("E46").Value = '  -10.52%  '

$ws.Range("D47").Value = '''42.77'
$ws.Range("E47").Value = '  -7.57%  '

$ws.Range("E48").Value = '  +5.90%  '

$ws.Range("E49").Value = '  -16.67%  '

$ws.Range("E50").Value = '  -13.95%  '

$ws.Range("D51").Value = '''0.576'
$ws.Range("E51").Value = '  -15.18%  '
